$d = $word.ActiveDocument

# --- (2b): "Fs > 60Hz  =  2.30 Hz" -> "Fs > 60Hz  =  2.(30) Hz" -------------
# Wrap "30" in parentheses to avoid reading "2.30" as a decimal number.
$rng = $d.Content
$rng.Find.Execute("2.30 Hz", $true, $false, $false, $false, $false, $true, 1, $false, "2.(30) Hz", 2)

# Word stamps the "_GoBack" bookmark at the location of the most recent edit;
# move it from its old spot (after the inline picture near the top of the
# doc) to right after the freshly inserted text, between "2." and "(30)".
$rng2 = $d.Content
$rng2.Find.Execute("2.(30) Hz", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackPos = $rng2.Start + 2
$bmRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
